$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date number format used by the existing "Date" column (column A, rows 3-18)
$dateFmt = $ws.Cells.Item(10, 1).NumberFormat

# --- Row 10: new entry for "Worked on enemy FSM" (short session) ---
$ws.Cells.Item(10, 1).Value2 = 44659
$ws.Cells.Item(10, 2).Value2 = "20.00 - 21.00"
$ws.Cells.Item(10, 4).Value2 = 1
$ws.Cells.Item(10, 5).Value2 = "Worked on enemy FSM"

# --- Row 11 ---
$ws.Cells.Item(11, 1).Value2 = 44660
$ws.Cells.Item(11, 2).Value2 = "10.00 - 12.00"
$ws.Cells.Item(11, 4).Value2 = 2
$ws.Cells.Item(11, 5).Value2 = "Worked on enemy FSM"

# --- Row 12 ---
$ws.Cells.Item(12, 1).Value2 = 44661
$ws.Cells.Item(12, 2).Value2 = "9.00 - 11.00"
$ws.Cells.Item(12, 4).Value2 = 2

# --- Row 13 ---
$ws.Cells.Item(13, 1).Value2 = 44662
$ws.Cells.Item(13, 2).Value2 = "9.30 - 12.00"
$ws.Cells.Item(13, 4).Value2 = 2.5

# --- Row 14 ---
$ws.Cells.Item(14, 1).Value2 = 44663

# --- Row 15 ---
$ws.Cells.Item(15, 1).Value2 = 44664

# --- Row 16 ---
$ws.Cells.Item(16, 1).Value2 = 44665

# --- Row 17 ---
$ws.Cells.Item(17, 1).Value2 = 44666

# --- Row 18 ---
$ws.Cells.Item(18, 1).Value2 = 44667
$ws.Cells.Item(18, 2).Value2 = "10.00 - 12.00"
$ws.Cells.Item(18, 4).Value2 = 2
$ws.Cells.Item(18, 5).Value2 = " "

# --- Row 19 (new row) ---
$ws.Cells.Item(19, 1).Value2 = 44669
$ws.Cells.Item(19, 1).NumberFormat = $dateFmt
$ws.Cells.Item(19, 2).Value2 = "9.30 - 12.00"
$ws.Cells.Item(19, 4).Value2 = 2.5

# --- Row 20 (new row) ---
$ws.Cells.Item(20, 1).Value2 = 44670
$ws.Cells.Item(20, 1).NumberFormat = $dateFmt
$ws.Cells.Item(20, 2).Value2 = "9.00 - 14.00"
$ws.Cells.Item(20, 4).Value2 = 5

# --- Row 21 (new row) ---
$ws.Cells.Item(21, 1).Value2 = 44671
$ws.Cells.Item(21, 1).NumberFormat = $dateFmt
$ws.Cells.Item(21, 2).Value2 = "9.00 - 14.00"
$ws.Cells.Item(21, 4).Value2 = 5

# --- Row 30: total formula recalculates automatically to 62.4 ---

# --- Selection moves to E18 ---
[void]$ws.Range("E18").Select()
